$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename date/time attribute names
$ws.Range("A2").Value = "datetime_utc"
$ws.Range("A3").Value = "datetime_utc_matlab"

# Insert two new rows for latitude/longitude after row 3 (before old row 4)
$ws.Rows("4:5").Insert()

$ws.Range("A4").Value = "latitude"
$ws.Range("B4").Value = "Latitude of sample event"
$ws.Range("C4").Value = "numeric"
$ws.Range("D4").Value = "degree"

$ws.Range("A5").Value = "longitude"
$ws.Range("B5").Value = "Longitude of sample event"
$ws.Range("C5").Value = "numeric"
$ws.Range("D5").Value = "degree"

# Update definition text cells (rows shifted down by 2 from originals)
$ws.Range("B6").Value = "Bottle sample from niskin or underway"
$ws.Range("B10").Value = "Oxygen-argon ratio divided by the reference ratio (oxygen-argon ratio in air minus 1, multiplied by 100)"
$ws.Range("B13").Value = "Oxygen-17 composition of dissolved oxygen versus atmospheric O2 in parts per thousand"
$ws.Range("B14").Value = "Oxygen-18 composition of dissolved oxygen versus atmospheric O2 in parts per thousand"

$ws.Range("E17").Select() | Out-Null
